# Apply the "tables_figures_structure" update:
#  1. Re-cache the datetimeFigureOut date placeholder (01/06/2021 -> 02/06/2021)
#     on the slide master and every slide layout.
#  2. Resize / reposition / re-fill / re-word several of the caption text
#     boxes on slide 1 to make room for the new "Commune-level top model"
#     captions.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder (type dt / 16) -- master + all custom layouts
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
        }
        if ($isDate -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "01/06/2021") {
                $sh.TextFrame.TextRange.Text = "02/06/2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($L).Shapes
}

# ---------------------------------------------------------------------
# 2. Slide 1 caption textboxes
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# -- "Figure - socioecon global model predictions with points" (shape 9)
#    grows taller (2 lines) and gets new wording for the 3rd run.
$sh9 = $s.Shapes.Item(9)
$sh9.Height = 523220 / 12700
$tr9 = $sh9.TextFrame.TextRange
$old9 = " global model predictions with points"
$new9 = " Commune – level top model predictions with points"
$start9 = $tr9.Text.IndexOf($old9) + 1
$tr9.Characters($start9, $old9.Length).Text = $new9

# -- "Figure - socioecon global model province-level predictions ..." (shape 10)
#    moves down, gets a green fill, and new wording for the 3rd run.
$sh10 = $s.Shapes.Item(10)
$sh10.Top = 3541580 / 12700
$sh10.Fill.Solid()
$sh10.Fill.ForeColor.RGB = 5296274  # 0x92D050 (R=146,G=208,B=80) as BGR COM long
$tr10 = $sh10.TextFrame.TextRange
$old10 = " global model province-level predictions with mean and commune-level lines"
$new10 = " commune-level top model province-level predictions with mean and commune-level lines"
$start10 = $tr10.Text.IndexOf($old10) + 1
$tr10.Characters($start10, $old10.Length).Text = $new10

# -- "Figure - Socioecon province-level categorical models" (shape 12)
#    moves down and gets new wording for the 3rd run.
$sh12 = $s.Shapes.Item(12)
$sh12.Top = 4115071 / 12700
$tr12 = $sh12.TextFrame.TextRange
$old12 = " province-level categorical models"
$new12 = " province-level categorical top model"
$start12 = $tr12.Text.IndexOf($old12) + 1
$tr12.Characters($start12, $old12.Length).Text = $new12

# -- "Figure - Cluster analysis map" (shape 14) -- moves down only
$sh14 = $s.Shapes.Item(14)
$sh14.Top = 4473119 / 12700

# -- "Figure - Cluster analysis heatmap" (shape 15) -- moves down only
$sh15 = $s.Shapes.Item(15)
$sh15.Top = 4831167 / 12700

# -- "Table (5) - Cluster analysis typologies" (shape 16) -- moves down only
$sh16 = $s.Shapes.Item(16)
$sh16.Top = 5167612 / 12700
